$d = $word.ActiveDocument

# 1) Replace the "Выводы" section body text with the new summary text.
$d.Content.Find.Execute(
    "Здесь кратко описываются итоги проделанной работы.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Мы освоили процедуры оформления отчетов с помощью легковесного языка разметки Markdown",
    2
) | Out-Null

# 2) Remove the entire "Список литературы" heading paragraph (its bookmarks
#    collapse down onto the end of the document, wrapping the now-empty
#    "refs" bookmark).
$last = $d.Paragraphs.Count
$p = $d.Paragraphs($last)
$p.Range.Delete()
